$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Z1").Value = "</mods:mods></datastream></object>"
$ws.Range("C1").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink">'

$ws.Range("C1").Select()
